$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the next empty row after the existing log data (row 9)
$newRow = 9

$ws.Cells.Item($newRow, 1).Value = "2025-08-13 13:08:48 UTC"
$ws.Cells.Item($newRow, 2).Value = "2025-08-13 18:38:48 IST"
$ws.Cells.Item($newRow, 3).Value = "SKIPPED"
$ws.Cells.Item($newRow, 4).Value = "No change in PDF. Skipping download & Excel update."
$ws.Cells.Item($newRow, 5).Value = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf"
$ws.Cells.Item($newRow, 6).Value = ""
$ws.Cells.Item($newRow, 7).Value = 0
$ws.Cells.Item($newRow, 8).Value = ""

# Match the style of the prior rows (centered alignment)
$srcRange = $ws.Range("A8:H8")
$dstRange = $ws.Range("A9:H9")
$srcRange.Copy()
$dstRange.PasteSpecial(-4122)
